# Update EC workbook: refresh "Estado de Cuenta" database with new mora data
# (part 1 of new estados de cuenta) and drop the JULIAN FELIPE CASTAÑO SALAZAR
# row, re-grouping the remaining three workers by period instead of by worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header summary figures -------------------------------------------------
$ws.Range("E11").Value = 2040120   # VALOR MORA total
$ws.Range("C13").Value = 3         # Cant. Trabajadores
$ws.Range("F13").Value = 6         # Cant. Periodos

# --- 2. Preserve the "last row" (bottom-border) look on the new final data row -
# Row 34 (JULIAN ...) currently carries the special bottom-border style used for
# the final row of the table. Copy that formatting onto row 33 (which will become
# the new final row) before row 34 is removed.
$ws.Range("B34:J34").Copy() | Out-Null
$ws.Range("B33:J33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 3. Remove the JULIAN FELIPE CASTAÑO SALAZAR row (row 34) ------------------
# Deleting this row also shifts the two footer rows (39, 40 -> 38, 39) up, which
# matches the target layout automatically.
$ws.Rows.Item(34).Delete()

# --- 4. Re-populate the data table (B16:J33), now grouped by period -----------
$workers = @(
    @("CC", "32706348", "DIANA LUZ DE LEON OBREGON", 124674, 3116850),
    @("CC", "22464377", "CARMEN ELENA MEZA ESTRADA", 124674, 3116850),
    @("CC", "9193183", "JOSE MARIA JIMENEZ MUNIVE", 90672, 2266800)
)
$periods = @("1607", "1608", "1609", "1610", "1611", "1612")

$r = 16
foreach ($periodo in $periods) {
    foreach ($w in $workers) {
        $ws.Cells.Item($r, 2).Value = $w[0]
        $ws.Cells.Item($r, 3).Value = $w[1]
        $ws.Cells.Item($r, 4).Value = $w[2]
        $ws.Cells.Item($r, 5).Value = $periodo
        $ws.Cells.Item($r, 6).Value = $w[3]
        $ws.Cells.Item($r, 7).Value = $w[4]
        $r = $r + 1
    }
}
